$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking values
# ("1.001", "29.349.30", etc.) are stored as text, not auto-converted numbers.
$ws.Range("D2:D16").NumberFormat = "@"
$ws.Range("D18:D27").NumberFormat = "@"
$ws.Range("D29:D34").NumberFormat = "@"
$ws.Range("D36:D45").NumberFormat = "@"
$ws.Range("D47:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.349.30"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "1.854.94"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "0.6987"
$ws.Range("E5").Value = "  -5.88%  "
$ws.Range("D6").Value = "239.61"
$ws.Range("E6").Value = "  -1.24%  "
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "0.3083"
$ws.Range("E8").Value = "  -2.29%  "
$ws.Range("D9").Value = "0.07496"
$ws.Range("E9").Value = "  +3.11%  "
$ws.Range("D10").Value = "23.78"
$ws.Range("E10").Value = "  -3.69%  "
$ws.Range("D11").Value = "0.08123"
$ws.Range("E11").Value = "  -3.29%  "
$ws.Range("D12").Value = "1.875.12"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "0.7285"
$ws.Range("E13").Value = "  -3.15%  "
$ws.Range("D14").Value = "5.219"
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").Value = "89.54"
$ws.Range("E15").Value = "  -3.33%  "
$ws.Range("D16").Value = "29.493.59"
$ws.Range("E16").Value = "  -1.50%  "
$ws.Range("E17").Value = "  -2.79%  "
$ws.Range("D18").Value = "243.19"
$ws.Range("E18").Value = "  -1.86%  "
$ws.Range("D19").Value = "0.000007748"
$ws.Range("E19").Value = "  -1.58%  "
$ws.Range("D20").Value = "13.16"
$ws.Range("E20").Value = "  -3.20%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("D22").Value = "2.132.81"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "7.643"
$ws.Range("E24").Value = "  -4.94%  "
$ws.Range("D25").Value = "0.1477"
$ws.Range("E25").Value = "  -5.46%  "
$ws.Range("D26").Value = "9.055"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").Value = "161.95"
$ws.Range("E27").Value = "  -1.97%  "
$ws.Range("E28").Value = "  -2.89%  "
$ws.Range("D29").Value = "1.945"
$ws.Range("E29").Value = "  -4.65%  "
$ws.Range("D30").Value = "1.390"
$ws.Range("E30").Value = "  -8.02%  "
$ws.Range("D31").Value = "1.509"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").Value = "4.418"
$ws.Range("E32").Value = "  -4.14%  "
$ws.Range("D33").Value = "4.063"
$ws.Range("E33").Value = "  -5.24%  "
$ws.Range("D34").Value = "0.05264"
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("E35").Value = "  -2.98%  "
$ws.Range("D36").Value = "0.7225"
$ws.Range("E36").Value = "  -4.17%  "
$ws.Range("D37").Value = "1.001"
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").Value = "2.667"
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("D39").Value = "0.01873"
$ws.Range("E39").Value = "  -4.64%  "
$ws.Range("D40").Value = "2.708"
$ws.Range("E40").Value = "  -1.83%  "
$ws.Range("D41").Value = "0.8849"
$ws.Range("E41").Value = "  +3.39%  "
$ws.Range("D42").Value = "0.4318"
$ws.Range("E42").Value = "  -4.64%  "
$ws.Range("D43").Value = "5.920"
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("D44").Value = "70.17"
$ws.Range("E44").Value = "  -3.47%  "
$ws.Range("D45").Value = "1.051.63"
$ws.Range("E45").Value = "  -5.48%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "102.84"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("D48").Value = "7.262"
$ws.Range("E48").Value = "  -4.80%  "
$ws.Range("D49").Value = "2.024.07"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "1.757"
$ws.Range("E50").Value = "  -5.54%  "
$ws.Range("D51").Value = "9.313"
$ws.Range("E51").Value = "  -1.98%  "
